{"js": "// Replace each old equation string with its corresponding new equation\n// string. Every \"old\" value is unique within the document, so an exact,\n// case-sensitive, non-wildcard search safely retargets exactly one run\n// per pair (order mirrors the document / diff order, though order does\n// not actually matter since each search is independent).\nconst replacements = [\n  [\"34\u00d738=1292\", \"63\u00d756=3528\"],\n  [\"39\u00d749=1911\", \"36\u00d785=3060\"],\n  [\"74\u00d788=6512\", \"94\u00d749=4606\"],\n  [\"84\u00d742=3528\", \"47\u00d730=1410\"],\n  [\"47\u00d794=4418\", \"14\u00d752=728\"],\n  [\"78\u00d753=4134\", \"26\u00d730=780\"],\n  [\"60\u00d750=3000\", \"27\u00d773=1971\"],\n  [\"19\u00d739=741\", \"33\u00d719=627\"],\n  [\"38\u00d761=2318\", \"62\u00d769=4278\"],\n  [\"65\u00d772=4680\", \"97\u00d749=4753\"],\n  [\"82\u00d784=6888\", \"73\u00d777=5621\"],\n  [\"38\u00d783=3154\", \"79\u00d775=5925\"],\n  [\"29\u00d738=1102\", \"19\u00d722=418\"],\n  [\"95\u00d791=8645\", \"67\u00d767=4489\"],\n  [\"74\u00d735=2590\", \"72\u00d712=864\"],\n  [\"15\u00d796=1440\", \"93\u00d772=6696\"],\n  [\"96\u00d786=8256\", \"98\u00d775=7350\"],\n  [\"72\u00d728=2016\", \"50\u00d790=4500\"],\n  [\"48\u00d798=4704\", \"20\u00d747=940\"],\n  [\"55\u00d778=4290\", \"95\u00d753=5035\"],\n  [\"31\u00d765=2015\", \"93\u00d750=4650\"],\n  [\"65\u00d760=3900\", \"75\u00d761=4575\"],\n  [\"80\u00d775=6000\", \"62\u00d751=3162\"],\n  [\"84\u00d773=6132\", \"55\u00d769=3795\"],\n  [\"22\u00d742=924\", \"66\u00d757=3762\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n    matchWildcards: false\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each \"old\" equation text is unique in the document, so a plain\n# (non-wildcard) Find/Replace across the whole document body safely\n# retargets exactly one run per pair, in document order.\n$pairs = @(\n    @(\"34\u00d738=1292\", \"63\u00d756=3528\"),\n    @(\"39\u00d749=1911\", \"36\u00d785=3060\"),\n    @(\"74\u00d788=6512\", \"94\u00d749=4606\"),\n    @(\"84\u00d742=3528\", \"47\u00d730=1410\"),\n    @(\"47\u00d794=4418\", \"14\u00d752=728\"),\n    @(\"78\u00d753=4134\", \"26\u00d730=780\"),\n    @(\"60\u00d750=3000\", \"27\u00d773=1971\"),\n    @(\"19\u00d739=741\", \"33\u00d719=627\"),\n    @(\"38\u00d761=2318\", \"62\u00d769=4278\"),\n    @(\"65\u00d772=4680\", \"97\u00d749=4753\"),\n    @(\"82\u00d784=6888\", \"73\u00d777=5621\"),\n    @(\"38\u00d783=3154\", \"79\u00d775=5925\"),\n    @(\"29\u00d738=1102\", \"19\u00d722=418\"),\n    @(\"95\u00d791=8645\", \"67\u00d767=4489\"),\n    @(\"74\u00d735=2590\", \"72\u00d712=864\"),\n    @(\"15\u00d796=1440\", \"93\u00d772=6696\"),\n    @(\"96\u00d786=8256\", \"98\u00d775=7350\"),\n    @(\"72\u00d728=2016\", \"50\u00d790=4500\"),\n    @(\"48\u00d798=4704\", \"20\u00d747=940\"),\n    @(\"55\u00d778=4290\", \"95\u00d753=5035\"),\n    @(\"31\u00d765=2015\", \"93\u00d750=4650\"),\n    @(\"65\u00d760=3900\", \"75\u00d761=4575\"),\n    @(\"80\u00d775=6000\", \"62\u00d751=3162\"),\n    @(\"84\u00d773=6132\", \"55\u00d769=3795\"),\n    @(\"22\u00d742=924\", \"66\u00d757=3762\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,    # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
